$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, per the diff.
$changes = @{
    "A3"  = -22.033
    "A21" = -19.937
    "A23" = -19.945
    "A25" = -21.789
    "C27" = -13.055
    "C31" = -13.002
    "C39" = -12.559
    "C48" = -11.1
    "C51" = -11.152
    "C52" = -11.601
    "A53" = -21.931
    "C55" = -13.513
    "C56" = -13.222
    "A57" = -22.095
    "C57" = -13.58
    "A59" = -22.407
    "A69" = -21.703
    "C73" = -12.515
    "A79" = -21.122
    "A83" = -21.938
    "C89" = -10.813
    "C90" = -12.482
    "A93" = -21.507
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
